# Adds *_ExternalURL values for the CAS "services" scenario data so that
# the CAS selection list logic can test for presence of a given service
# before listing a CAS for that service.
#
# Columns (row 1 headers):
#   K = OWA_InternalURL   L = OWA_ExternalURL
#   N = EAS_InternalURL   O = EAS_ExternalURL
#   P = EWS_InternalURL   Q = EWS_ExternalURL
#   R = ECP_InternalURL   S = ECP_ExternalURL
#
# Only the CAS rows that already publish the corresponding *_InternalURL
# get a matching *_ExternalURL value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAS_scenario")

# OWA_ExternalURL (column L) - one per CAS that publishes OWA_InternalURL
$ws.Range("L2").Value  = "https://us.mail.corp.com/owa"        # WAS10CAS01
$ws.Range("L4").Value  = "https://us.legacy.corp.com/owa"      # WAS07CAS01
$ws.Range("L14").Value = "https://eu.mail.corp.com/owa"        # LON10CAS01
$ws.Range("L16").Value = "https://eu.legacy.corp.com/owa"      # LON07CAS01
$ws.Range("L20").Value = "https://lgb.mail.corp.com/owa"       # LGB10EXC01

# EAS_ExternalURL (column O) - one per CAS that publishes EAS_InternalURL
$ws.Range("O2").Value  = "https://us.mail.corp.com/Microsoft-Server-ActiveSync"
$ws.Range("O4").Value  = "https://us.legacy.corp.com/Microsoft-Server-ActiveSync"
$ws.Range("O14").Value = "https://eu.mail.corp.com/Microsoft-Server-ActiveSync"
$ws.Range("O16").Value = "https://eu.legacy.corp.com/Microsoft-Server-ActiveSync"
$ws.Range("O20").Value = "https://lgb.mail.corp.com/Microsoft-Server-ActiveSync"

# EWS_ExternalURL (column Q) - one per CAS that publishes EWS_InternalURL
$ws.Range("Q2").Value  = "https://us.mail.corp.com/EWS/Exchange.asmx"
$ws.Range("Q4").Value  = "https://us.legacy.corp.com/EWS/Exchange.asmx"
$ws.Range("Q14").Value = "https://eu.mail.corp.com/EWS/Exchange.asmx"
$ws.Range("Q16").Value = "https://eu.legacy.corp.com/EWS/Exchange.asmx"
$ws.Range("Q20").Value = "https://lgb.mail.corp.com/EWS/Exchange.asmx"

# ECP_ExternalURL (column S) - one per CAS that publishes ECP_InternalURL
$ws.Range("S2").Value  = "https://us.mail.corp.com/ecp"
$ws.Range("S14").Value = "https://eu.mail.corp.com/ecp"
$ws.Range("S20").Value = "https://lgb.mail.corp.com/ecp"
